# Edit "Liste des taches.xlsx" per the target commit:
#  - Add mailchimp address + password notes next to row 24
#  - Move the status mark for "Mettre à jour réglements et remboursements"
#    (row 35) from "En cours" (col B) to "Mise en forme restante" (col D),
#    and fill in its start/end dates (cols F/G)
#  - Update the active selection to J20
#  - Resize the saved window (cosmetic)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New documentation cells on row 24 (next to "Mise en forme restante" header block)
$ws.Range("G24").Value = "Adresse mailchimp : delasalle.sio.destinataire@gmail.com"
$ws.Range("J24").Value = "mdp : Sio1_Sio2"

# Row 35 ("      Mettre à jour réglements et remboursements"): task moved
# from "En cours" to "Mise en forme restante"
$ws.Range("B35").Value = $null
$ws.Range("D35").Value = "X"

# Fill in start/end dates for that row, matching the date-formatted style
# used by the other rows (copy number format from row 32's F:G cells).
$ws.Range("F35").Value = 42521
$ws.Range("G35").Value = 42522
$ws.Range("F32:G32").Copy()
$ws.Range("F35:G35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Window/selection state as left by the author
$win = $excel.ActiveWindow
$win.Width = 23145
$win.Height = 9210

$ws.Range("J20").Select() | Out-Null
